$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row: new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (e.g. AC1) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1").Select()

# Fill in team record data for each data row (2 through 56)
for ($r = 2; $r -le 56; $r++) {
    $ws.Cells.Item($r, 30).Value = 73   # AD
    $ws.Cells.Item($r, 31).Value = 89   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
